{"js": "// Replace each three-digit-by-one-digit multiplication prompt with its\n// updated version. Every \"old\" equation text below is unique within the\n// document, so a simple whole-body search + replace is safe and order\n// independent.\nconst replacements = [\n  [\"736\u00d74=\", \"644\u00d75=\"],\n  [\"676\u00d74=\", \"915\u00d78=\"],\n  [\"956\u00d75=\", \"647\u00d75=\"],\n  [\"964\u00d75=\", \"307\u00d78=\"],\n  [\"355\u00d79=\", \"503\u00d78=\"],\n  [\"723\u00d77=\", \"883\u00d78=\"],\n  [\"616\u00d79=\", \"406\u00d79=\"],\n  [\"353\u00d77=\", \"417\u00d75=\"],\n  [\"644\u00d77=\", \"465\u00d77=\"],\n  [\"856\u00d76=\", \"907\u00d79=\"],\n  [\"890\u00d72=\", \"229\u00d73=\"],\n  [\"497\u00d76=\", \"822\u00d75=\"],\n  [\"578\u00d76=\", \"304\u00d75=\"],\n  [\"780\u00d73=\", \"299\u00d73=\"],\n  [\"457\u00d75=\", \"899\u00d75=\"],\n  [\"191\u00d72=\", \"547\u00d74=\"],\n  [\"276\u00d74=\", \"803\u00d73=\"],\n  [\"410\u00d76=\", \"557\u00d76=\"],\n  [\"306\u00d79=\", \"569\u00d78=\"],\n  [\"748\u00d72=\", \"534\u00d78=\"],\n  [\"759\u00d77=\", \"380\u00d74=\"],\n  [\"551\u00d74=\", \"488\u00d77=\"],\n  [\"548\u00d73=\", \"390\u00d74=\"],\n  [\"233\u00d72=\", \"453\u00d78=\"],\n  [\"963\u00d78=\", \"252\u00d79=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update each three-digit-by-one-digit multiplication prompt in place.\n# Every \"old\" equation text is unique within the document, so a plain\n# Find/Replace (no wildcards) for each pair is safe and order independent.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"736\u00d74=\", \"644\u00d75=\"),\n    @(\"676\u00d74=\", \"915\u00d78=\"),\n    @(\"956\u00d75=\", \"647\u00d75=\"),\n    @(\"964\u00d75=\", \"307\u00d78=\"),\n    @(\"355\u00d79=\", \"503\u00d78=\"),\n    @(\"723\u00d77=\", \"883\u00d78=\"),\n    @(\"616\u00d79=\", \"406\u00d79=\"),\n    @(\"353\u00d77=\", \"417\u00d75=\"),\n    @(\"644\u00d77=\", \"465\u00d77=\"),\n    @(\"856\u00d76=\", \"907\u00d79=\"),\n    @(\"890\u00d72=\", \"229\u00d73=\"),\n    @(\"497\u00d76=\", \"822\u00d75=\"),\n    @(\"578\u00d76=\", \"304\u00d75=\"),\n    @(\"780\u00d73=\", \"299\u00d73=\"),\n    @(\"457\u00d75=\", \"899\u00d75=\"),\n    @(\"191\u00d72=\", \"547\u00d74=\"),\n    @(\"276\u00d74=\", \"803\u00d73=\"),\n    @(\"410\u00d76=\", \"557\u00d76=\"),\n    @(\"306\u00d79=\", \"569\u00d78=\"),\n    @(\"748\u00d72=\", \"534\u00d78=\"),\n    @(\"759\u00d77=\", \"380\u00d74=\"),\n    @(\"551\u00d74=\", \"488\u00d77=\"),\n    @(\"548\u00d73=\", \"390\u00d74=\"),\n    @(\"233\u00d72=\", \"453\u00d78=\"),\n    @(\"963\u00d78=\", \"252\u00d79=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Execute(\n        $oldText,\n        $false,\n        $true,\n        $false,\n        $false,\n        $false,\n        $true,\n        1,\n        $false,\n        $newText,\n        2\n    )\n}\n"}
